$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Fill in the rest of April-25 (12)'s daily purchase figures
# ------------------------------------------------------------------
$april = $wb.Worksheets.Item("April-25 (12)")

$april.Range("R7").Value = 134300

$april.Range("C8").Value = 200000
$april.Range("N8").Value = 10
$april.Range("O8").Value = 100
$april.Range("R8").Value = 160651

$april.Range("C9").Value = 92212
$april.Range("I9").Value = 2000
$april.Range("L9").Value = 500
$april.Range("R9").Value = 102006

$april.Range("C10").Value = 114285
$april.Range("R10").Value = 103873

$april.Range("C11").Value = 114285
$april.Range("R11").Value = 113876

$april.Range("C12").Value = 231901
$april.Range("I12").Value = 2000
$april.Range("L12").Value = 500
$april.Range("N12").Value = 20
$april.Range("O12").Value = 50
$april.Range("R12").Value = 141117

$april.Range("R14").Value = 99989

$april.Range("C15").Value = 195766
$april.Range("R15").Value = 107067

$april.Range("C16").Value = 31240
$april.Range("R16").Value = 103553

$april.Range("C17").Value = 81018
$april.Range("R17").Value = 95978

$april.Range("C18").Value = 146909
$april.Range("N18").Value = 35
$april.Range("O18").Value = 50
$april.Range("R18").Value = 196462

$april.Range("C19").Value = 265690
$april.Range("I19").Value = 2500
$april.Range("K19").Value = 500
$april.Range("L19").Value = 500
$april.Range("R19").Value = 144864

$april.Range("R21").Value = 89838

$april.Range("C22").Value = 88415
$april.Range("I22").Value = 2500
$april.Range("N22").Value = 25
$april.Range("O22").Value = 25
$april.Range("R22").Value = 109769

$april.Range("C23").Value = 103896
$april.Range("R23").Value = 107970

$april.Range("C24").Value = 114285
$april.Range("R24").Value = 111138

$april.Range("C25").Value = 124675
$april.Range("R25").Value = 126209

$april.Range("C26").Value = 270017
$april.Range("I26").Value = 5000
$april.Range("K26").Value = 500
$april.Range("L26").Value = 500
$april.Range("R26").Value = 160545

$april.Range("R28").Value = 103769

$april.Range("C29").Value = 103896
$april.Range("R29").Value = 107923

$april.Range("C30").Value = 109964
$april.Range("I30").Value = 2000
$april.Range("N30").Value = 25
$april.Range("R30").Value = 118406

$april.Range("C31").Value = 169870
$april.Range("R31").Value = 174076

$april.Range("C32").Value = 402743
$april.Range("I32").Value = 2000
$april.Range("N32").Value = 25
$april.Range("R32").Value = 207552

# Update April-25's view state (no longer the active/selected tab; zoomed
# back to 100% with the selection moved to T32)
$april.Select()
$april.Range("T32").Select()
$excel.ActiveWindow.Zoom = 100

# ------------------------------------------------------------------
# 2) Add the new May-25 (13) sheet, copied from April-25's layout, and
#    fill the first few days of data that had already come in.
# ------------------------------------------------------------------
$may = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $april)
$may.Name = "May-25 (13)"

# Header / title row
$may.Range("A1").Value = $april.Range("A1").Value
$may.Range("A1:R1").Merge()

# Column headers (row 2)
for ($col = 1; $col -le 18; $col++) {
    $may.Cells.Item(2, $col).Value = $april.Cells.Item(2, $col).Value
}

# Day label + date columns, 31 days (May has 31 days) -> rows 3..33
$startSerial = 45413
for ($i = 0; $i -lt 31; $i++) {
    $r = 3 + $i
    $may.Cells.Item($r, 1).Value = $april.Cells.Item(3 + ($i % 30), 1).Value
    $may.Cells.Item($r, 2).Value = $startSerial + $i
    $may.Cells.Item($r, 2).NumberFormat = $april.Cells.Item(3, 2).NumberFormat
}

# Totals row (34) - shared SUM formula across C..R, same as April sheet
$may.Range("A34").Value = $april.Range("A34").Value
$may.Range("B34").Value = $april.Range("B34").Value
$may.Range("C34:R34").FormulaR1C1 = "=SUM(R[-31]C:R[-1]C)"

# Second header row (35) mirrors row 2
for ($col = 1; $col -le 18; $col++) {
    $may.Cells.Item(35, $col).Value = $april.Cells.Item(35, $col).Value
}

# CARD PURCHASE / TARGAT / DUE block (rows 36-40)
$may.Range("C36").Value = $april.Range("C36").Value
$may.Range("E36").Formula = $april.Range("E36").Formula
$may.Range("C38").Value = $april.Range("C38").Value
$may.Range("E38").Value = $april.Range("E38").Value
$may.Range("E39").Value = $april.Range("E39").Value
$may.Range("C40").Value = $april.Range("C40").Value
$may.Range("E40").Formula = $april.Range("E40").Formula

# Data already reported for the first few days of May
$may.Range("C5").Value = 5236
$may.Range("R5").Value = 129016
$may.Range("R6").Value = 120964
$may.Range("R7").Value = 109190
$may.Range("C8").Value = 124675
$may.Range("R8").Value = 114495

# May-25 becomes the active/selected sheet
$may.Select()
$may.Range("C9").Select()
$excel.ActiveWindow.Zoom = 100
